$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.552.05"
$ws.Range("E2").Value = "  +5.18%  "

$ws.Range("D3").Value = "3.645.80"
$ws.Range("E3").Value = "  +5.56%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.02"
$ws.Range("E5").Value = "  +2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.32"
$ws.Range("E6").Value = "  +4.12%  "

$ws.Range("E7").Value = "  +2.15%  "

$ws.Range("D8").Value = "3.642.72"
$ws.Range("E8").Value = "  +5.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  +6.22%  "

$ws.Range("E11").Value = "  +4.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.04"
$ws.Range("E12").Value = "  +3.40%  "

$ws.Range("E13").Value = "  +5.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.94"
$ws.Range("E14").Value = "  +5.71%  "

$ws.Range("D15").Value = "4.228.36"
$ws.Range("E15").Value = "  +5.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.01"
$ws.Range("E16").Value = "  +6.96%  "

$ws.Range("D17").Value = "3.645.17"
$ws.Range("E17").Value = "  +5.67%  "

$ws.Range("D18").Value = "70.554.56"
$ws.Range("E18").Value = "  +5.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +5.23%  "

$ws.Range("E20").Value = "  +3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  +3.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.28"
$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.17"
$ws.Range("E23").Value = "  +13.45%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.53"
$ws.Range("E25").Value = "  +3.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.16"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("E27").Value = "  +7.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.49"
$ws.Range("E28").Value = "  +5.02%  "

$ws.Range("E29").Value = "  +5.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.73"
$ws.Range("E30").Value = "  +4.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("E31").Value = "  +7.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  +9.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "626.38"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.24"
$ws.Range("E34").Value = "  +4.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.85"
$ws.Range("E35").Value = "  +3.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.60"
$ws.Range("E36").Value = "  +10.69%  "

$ws.Range("E37").Value = "  +9.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.411"
$ws.Range("E38").Value = "  +6.77%  "

$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("D42").Value = "3.301.48"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("E43").Value = "  +7.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +12.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0454"
$ws.Range("E45").Value = "  +5.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +9.05%  "

$ws.Range("E47").Value = "  +2.32%  "

$ws.Range("E48").Value = "  +2.99%  "

$ws.Range("E49").Value = "  +6.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +0.08%  "
